# The 06.07.19 s1cDNA sample with s1cDNASampleNumber 12 (row 13) had no
# fastq file in the 06.17.19 library, so it is removed from the sheet.
# Deleting the entire row shifts every subsequent sample row up by one,
# which is exactly the row-by-row renumbering seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 13 first (mirrors the author selecting the row before
# deleting it - this is also what leaves the sheet's cached selection
# at A13:XFD13 after the edit).
$ws.Rows(13).Select() | Out-Null

# Delete the entire row; cells below shift up to fill the gap.
$ws.Rows(13).Delete()
